$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.702.53"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "3.306.50"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.308.06"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.32%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "3.851.27"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "3.304.16"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "63.826.01"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.737"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.27%  "
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "0.0₃0740"
$ws.Range("E38").Value = "  +3.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0400"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").Value = "3.120.84"
$ws.Range("E40").Value = "  +4.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "430.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("E42").Value = "  +7.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "128.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.21%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
